try {
    1/0
} catch {
    Write-Host "Caught div by zero:" $_.Exception.Message
}
Write-Host "After try"
